$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in cell E8 (was "Good Morning", now "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Select the edited cell, matching the active selection recorded in the file
$ws.Activate()
$ws.Range("E8").Select()
